$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the typo in the existing antiarrhythmic-agent-therapy row (row 63) ---
$ws.Range("B63").Value = "antiarrythmic agent therapy"

# --- Append the new calcium-channel-modulator-therapy rows (64-73) ---
# Each row: A=defined_class, B=defined_class_name, C=chemical_identity, D=chemical_identity_label
$rows = @(
    @{ r = 64;  a = "MAXO_0000650"; b = "antiatherogenic agent therapy";       c = "CHEBI:50855";  d = "antiatherogenic agent";            tall = $true  },
    @{ r = 65;  a = "MAXO_0000651"; b = "antiatherosclerotic agent therapy";   c = "CHEBI:145947"; d = "antiaterosclerotic agent";          tall = $true  },
    @{ r = 66;  a = "MAXO_0001025"; b = "adrenergic agent therapy";           c = "CHEBI:37962";  d = "adrenergic agent";                  tall = $true  },
    @{ r = 67;  a = "MAXO_0000182"; b = "alpha adrenergic agent therapy";     c = "CHEBI:48539";  d = "alpha-adrenergic drug";             tall = $true  },
    @{ r = 68;  a = "MAXO_0000183"; b = "alpha adrenergic agonist therapy";   c = "CHEBI:48539";  d = "alpha-adrenergic agonist drug";     tall = $false },
    @{ r = 69;  a = "MAXO_0000184"; b = "alpha adrenergic antagonist therapy";c = "CHEBI:37890";  d = "alpha-adrenergic antagonist drug";  tall = $false },
    @{ r = 70;  a = "MAXO_0000186"; b = "beta adrenergic agent therapy";      c = "CHEBI:48540";  d = "beta-adrenergic agent thearpy";     tall = $false },
    @{ r = 71;  a = "MAXO_0001026"; b = "beta-adrenergic agonist therapy";    c = "CHEBI:35522";  d = "beta-adrenergic agonist drug";      tall = $false },
    @{ r = 72;  a = "MAXO_0000187"; b = "beta-adrenergic antagonist therapy"; c = "CHEBI:35530";  d = "beta-adrenergic antagonist drug";   tall = $false },
    @{ r = 73;  a = "MAXO_0000434"; b = "calcium channel blocking agent";     c = "CHEBI:38215";  d = "calcium channel blocker";           tall = $true  }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = $row.d

    if ($row.tall) {
        # Matches the existing "chemical_identity" styling used elsewhere in the
        # sheet (Helvetica 13pt) which the taller rows in this table use.
        $ws.Range("C$r").Font.Name = "Helvetica"
        $ws.Range("C$r").Font.Size = 13
        $ws.Rows.Item($r).RowHeight = 17
    }
}

# --- Update the view state to match where the author ended up after editing ---
$ws.Range("B74").Select()
